$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1831
$ws.Range("J2").Value = 1831
$ws.Range("L2").Value = 1831
$ws.Range("N2").Value = -2057
$ws.Range("H6").Value = 625.7273
$ws.Range("I6").Value = 405
$ws.Range("K6").Value = 1215
$ws.Range("M6").Value = -1103
$ws.Range("H40").Value = 26086.75
$ws.Range("I40").Value = 50447.5
$ws.Range("K40").Value = 50447.5
$ws.Range("M40").Value = -50272.5
$ws.Range("H55").Value = 546.125
$ws.Range("I55").Value = 473.4
$ws.Range("J55").Value = 667.3333
$ws.Range("K55").Value = 473.4
$ws.Range("L55").Value = 667.3333
$ws.Range("M55").Value = -259.4
$ws.Range("N55").Value = -1095.3333
$ws.Range("H97").Value = 2144.1667
$ws.Range("J97").Value = 2144.1667
$ws.Range("L97").Value = 6432.500100000001
$ws.Range("N97").Value = -7424.500100000001
$ws.Range("H98").Value = 641.125
$ws.Range("I98").Value = 634.2174
$ws.Range("K98").Value = 634.2174
$ws.Range("M98").Value = 863.7826
$ws.Range("H106").Value = 2289.3845
$ws.Range("I106").Value = 1976.2
$ws.Range("K106").Value = 1976.2
$ws.Range("M106").Value = -1345.2
$ws.Range("H122").Value = 641.125
$ws.Range("I122").Value = 634.2174
$ws.Range("K122").Value = 1902.6522
$ws.Range("M122").Value = 547.3478
$ws.Range("H125").Value = 1625.4
$ws.Range("I125").Value = 1352.5385
$ws.Range("K125").Value = 12172.8465
$ws.Range("M125").Value = -9712.846500000001
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 2231.92
$ws.Range("I129").Value = 626.7692
$ws.Range("J129").Value = 3970.8333
$ws.Range("K129").Value = 1880.3076
$ws.Range("L129").Value = 11912.4999
$ws.Range("M129").Value = 3119.6924
$ws.Range("N129").Value = -21912.4999
$ws.Range("H132").Value = 11930.223
$ws.Range("I132").Value = 3218.3684
$ws.Range("J132").Value = 32620.875
$ws.Range("K132").Value = 9655.1052
$ws.Range("L132").Value = 97862.625
$ws.Range("M132").Value = -7125.1052
$ws.Range("N132").Value = -102922.625
$ws.Range("H137").Value = 14494819
$ws.Range("I137").Value = 1828
$ws.Range("K137").Value = 5484
$ws.Range("M137").Value = -2934

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16806.17
$ws.Range("I32").Value = 16244.276
$ws.Range("K32").Value = 16244.276
$ws.Range("M32").Value = -15957.276
$ws.Range("H45").Value = 3142.8333
$ws.Range("I45").Value = 2460.6667
$ws.Range("K45").Value = 2460.6667
$ws.Range("M45").Value = -2083.6667
$ws.Range("H46").Value = 24637.166
$ws.Range("J46").Value = 28004.6
$ws.Range("L46").Value = 28004.6
$ws.Range("N46").Value = -28642.6
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 1747.0834
$ws.Range("I74").Value = 1385.1111
$ws.Range("J74").Value = 2833
$ws.Range("K74").Value = 1385.1111
$ws.Range("L74").Value = 2833
$ws.Range("M74").Value = -511.1111000000001
$ws.Range("N74").Value = -4581
$ws.Range("H77").Value = 1747.0834
$ws.Range("I77").Value = 1385.1111
$ws.Range("J77").Value = 2833
$ws.Range("K77").Value = 6925.5555
$ws.Range("L77").Value = 14165
$ws.Range("M77").Value = -2557.5555
$ws.Range("N77").Value = -22901
$ws.Range("H102").Value = 1802.3334
$ws.Range("I102").Value = 1768
$ws.Range("K102").Value = 1768
$ws.Range("M102").Value = -146
$ws.Range("H122").Value = 8006.4165
$ws.Range("I122").Value = 5594.1665
$ws.Range("K122").Value = 16782.4995
$ws.Range("M122").Value = -14332.4995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 38997
$ws.Range("I96").Value = 35796.4
$ws.Range("K96").Value = 35796.4
$ws.Range("M96").Value = -33050.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41669356
$ws.Range("I31").Value = 55557150
$ws.Range("J31").Value = 5960.3335
$ws.Range("K31").Value = 55557150
$ws.Range("L31").Value = 5960.3335
$ws.Range("M31").Value = -55556855
$ws.Range("N31").Value = -6550.3335
$ws.Range("H34").Value = 41669356
$ws.Range("I34").Value = 55557150
$ws.Range("J34").Value = 5960.3335
$ws.Range("K34").Value = 55557150
$ws.Range("L34").Value = 5960.3335
$ws.Range("M34").Value = -55556948
$ws.Range("N34").Value = -6364.3335
$ws.Range("H58").Value = 2354.0667
$ws.Range("I58").Value = 2143.875
$ws.Range("K58").Value = 2143.875
$ws.Range("M58").Value = -1940.875
$ws.Range("H86").Value = 7873.5
$ws.Range("J86").Value = 7938.2856
$ws.Range("L86").Value = 7938.2856
$ws.Range("N86").Value = -10184.2856
$ws.Range("H89").Value = 7873.5
$ws.Range("J89").Value = 7938.2856
$ws.Range("L89").Value = 39691.428
$ws.Range("N89").Value = -50923.428
$ws.Range("H100").Value = 78354.8
$ws.Range("J100").Value = 78354.8
$ws.Range("L100").Value = 78354.8
$ws.Range("N100").Value = -80518.8
$ws.Range("H132").Value = 34199612
$ws.Range("J132").Value = 30002.545
$ws.Range("L132").Value = 90007.63499999999
$ws.Range("N132").Value = -95067.63499999999
$ws.Range("H136").Value = 2354.0667
$ws.Range("I136").Value = 2143.875
$ws.Range("K136").Value = 6431.625
$ws.Range("M136").Value = -3881.625
$ws.Range("H141").Value = 106718.05
$ws.Range("J141").Value = 108989.414
$ws.Range("L141").Value = 108989.414
$ws.Range("N141").Value = -119349.414

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1792.5
$ws.Range("I48").Value = 1790
$ws.Range("J48").Value = 1795
$ws.Range("K48").Value = 5370
$ws.Range("L48").Value = 5385
$ws.Range("M48").Value = -5120
$ws.Range("N48").Value = -5885
$ws.Range("H56").Value = 55564900
$ws.Range("I56").Value = 55564900
$ws.Range("K56").Value = 55564900
$ws.Range("M56").Value = -55564370
$ws.Range("H107").Value = 1236.5
$ws.Range("J107").Value = 1236.5
$ws.Range("L107").Value = 3709.5
$ws.Range("N107").Value = -7549.5
$ws.Range("H137").Value = 70834456
$ws.Range("I137").Value = 68183050
$ws.Range("J137").Value = 100000000
$ws.Range("K137").Value = 204549150
$ws.Range("L137").Value = 300000000
$ws.Range("M137").Value = -204544050
$ws.Range("N137").Value = -300010200

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 22539.6
$ws.Range("J14").Value = 3099.5
$ws.Range("L14").Value = 3099.5
$ws.Range("N14").Value = -3435.5
$ws.Range("H53").Value = 3000
$ws.Range("I53").Value = 3000
$ws.Range("K53").Value = 3000
$ws.Range("M53").Value = -2369
$ws.Range("H102").Value = 9768.048000000001
$ws.Range("I102").Value = 10010.667
$ws.Range("J102").Value = 9161.5
$ws.Range("K102").Value = 10010.667
$ws.Range("L102").Value = 9161.5
$ws.Range("M102").Value = -8388.666999999999
$ws.Range("N102").Value = -12405.5
$ws.Range("H107").Value = 1018.125
$ws.Range("J107").Value = 1035
$ws.Range("L107").Value = 1035
$ws.Range("N107").Value = -4875
$ws.Range("H122").Value = 345691.6
$ws.Range("I122").Value = 501117.88
$ws.Range("J122").Value = 3753.8
$ws.Range("K122").Value = 1503353.64
$ws.Range("L122").Value = 11261.4
$ws.Range("M122").Value = -1500903.64
$ws.Range("N122").Value = -16161.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6583019
$ws.Range("I40").Value = 8336372.5
$ws.Range("K40").Value = 8336372.5
$ws.Range("M40").Value = -8336236.5
$ws.Range("H95").Value = 59500
$ws.Range("J95").Value = 59500
$ws.Range("L95").Value = 59500
$ws.Range("N95").Value = -64992
$ws.Range("H122").Value = 5264.95
$ws.Range("I122").Value = 3167
$ws.Range("K122").Value = 9501
$ws.Range("M122").Value = -7051

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14794.333
$ws.Range("J81").Value = 17692.727
$ws.Range("L81").Value = 35385.454
$ws.Range("N81").Value = -37507.454
$ws.Range("H84").Value = 14794.333
$ws.Range("J84").Value = 17692.727
$ws.Range("L84").Value = 176927.27
$ws.Range("N84").Value = -187535.27
$ws.Range("H100").Value = 643
$ws.Range("I100").Value = 199
$ws.Range("K100").Value = 398
$ws.Range("M100").Value = 143
$ws.Range("H122").Value = 5240.7896
$ws.Range("I122").Value = 3631.3333
$ws.Range("K122").Value = 10893.9999
$ws.Range("M122").Value = -8443.999899999999
$ws.Range("H136").Value = 9917.825999999999
$ws.Range("I136").Value = 2988.3333
$ws.Range("K136").Value = 8964.999899999999
$ws.Range("M136").Value = -6414.999899999999
